# Case_2_177: bus voltage magnitude (vm_pu) results updated for the 380 kV case.
# Column B holds the slack-bus setpoint (1.05 -> 1.02 p.u.); columns C:F and I:N
# are the recomputed load-flow results for buses 1-4 and 7-12 on every row (2-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns touched by this case re-run, in sheet order (B,C,D,E,F,I,J,K,L,M,N).
$columns = @(2,3,4,5,6,9,10,11,12,13,14)

# One array per data row (r=2..25): new vm_pu values aligned with $columns.
$newValues = @{
    2 = @(1.02,1.042034035162742,1.049615836873899,1.055701767183661,1.062164314595168,1.039235777331006,1.047112032540698,1.052372400553598,1.058441500330583,1.064886404280375,1.019522391235485)
    3 = @(1.02,1.043301195006434,1.050623411325124,1.056842726835421,1.063349361268014,1.039536506873624,1.048023885302263,1.053191752598147,1.059395118717814,1.065885282122408,1.019830818324642)
    4 = @(1.02,1.044120452674223,1.051274517116386,1.057581279893777,1.064116025035685,1.03972900360946,1.04861275175547,1.053720458081628,1.060011865547078,1.066530906554367,1.020029846752443)
    5 = @(1.02,1.044464708796184,1.051548036951635,1.057891835036667,1.06443829846908,1.039809428103887,1.048860034429492,1.053942375253259,1.060271073462354,1.066802157060885,1.020113388445093)
    6 = @(1.02,1.044522501626172,1.051593950151932,1.057943982606582,1.064492407767397,1.039822902358243,1.048901538097167,1.053979615601003,1.060314591400718,1.066847691288296,1.020127407862446)
    7 = @(1.02,1.044125053265726,1.051278172705778,1.057585429280087,1.064120331392816,1.039730080213994,1.048616057044133,1.053723424725787,1.060015329379016,1.066534531683084,1.020030963550264)
    8 = @(1.02,1.042462420186234,1.049956530750894,1.056087303560449,1.062564836797749,1.039337844459627,1.047420439339672,1.052649609604799,1.058763844690088,1.065224128940884,1.01962673866429)
    9 = @(1.02,1.039527315149898,1.047620961842643,1.053449464379275,1.059822703908771,1.038630603390936,1.045304621961877,1.050746090803719,1.056556151083368,1.062909478363558,1.01891024964365)
    10 = @(1.02,1.037566809604749,1.046059341806343,1.051692191946843,1.057993731823832,1.038148268965417,1.043887924707701,1.04946938038763,1.055082645391564,1.061362546960202,1.018429738264194)
    11 = @(1.02,1.036716953557874,1.045382036734046,1.050931557679214,1.057201532033999,1.03793683325294,1.043272994461362,1.048914704462552,1.054444175981352,1.06069177737318,1.018220987444681)
    12 = @(1.02,1.036401133637749,1.04513028624411,1.050649063831705,1.056907235002395,1.037857907928031,1.043044355684686,1.048708393028345,1.054206953369143,1.06044248069337,1.018143344299998)
    13 = @(1.02,1.03646888472831,1.045184295268583,1.050709657960575,1.056970364505274,1.037874855280206,1.043093409756218,1.048752660236633,1.054257841469351,1.060495962181139,1.018160003726344)
    14 = @(1.02,1.03669085077522,1.045361230423695,1.050908205853265,1.057177206153425,1.037930317194798,1.043254099738639,1.048897656419866,1.054424568456133,1.060671173360726,1.01821457155958)
    15 = @(1.02,1.036827591985821,1.045470223528628,1.051030542950746,1.057304642915393,1.03796443758156,1.043353076083781,1.048986956181717,1.054527285570642,1.060779107742034,1.018248178827683)
    16 = @(1.02,1.037623191543928,1.046104268689873,1.051742678395399,1.058046302205105,1.038162246794257,1.043928703993646,1.049506153225302,1.055125009255676,1.061407043789385,1.018443577851162)
    17 = @(1.02,1.03812199473784,1.046501689023286,1.052189454762439,1.058511458873096,1.038285635516154,1.044289379375374,1.049831334692632,1.055499828691113,1.061800678885968,1.018565962283394)
    18 = @(1.02,1.038412847350588,1.046733390298701,1.052450078648388,1.058782753616278,1.038357357027253,1.044499611522023,1.050020828988841,1.055718412993628,1.062030189207429,1.018637280889568)
    19 = @(1.02,1.038512005276103,1.04681237630799,1.052538949259917,1.058875254233744,1.038381769991285,1.044571270945812,1.05008541141735,1.055792937600766,1.062108431028975,1.018661587489825)
    20 = @(1.02,1.038068487300955,1.046459060694542,1.052141517130448,1.058461554389506,1.038272422840061,1.044250697201665,1.049796464281631,1.055459618422709,1.061758454928759,1.018552838440958)
    21 = @(1.02,1.036625491354217,1.045309132117968,1.050849737370263,1.057116297531695,1.037913995779544,1.043206786809538,1.04885496640506,1.054375473368212,1.060619581999128,1.018198505563664)
    22 = @(1.02,1.035717378614646,1.044585146285529,1.050037772239745,1.056270255776531,1.037686389227229,1.042549128441239,1.048261387817873,1.053693441861198,1.05990269906754,1.017975121540945)
    23 = @(1.02,1.03619886758888,1.044969038437028,1.050468189344033,1.056718780645599,1.037807261241976,1.042897890574513,1.048576209386835,1.054055036899612,1.060282811477319,1.018093598828694)
    24 = @(1.02,1.038092665289079,1.046478322932894,1.052163178000168,1.058484104149455,1.038278393849168,1.044268176454889,1.049812221271078,1.055477787840529,1.061777534396755,1.01855876874488)
    25 = @(1.02,1.040286760954974,1.048225561823064,1.054131175469731,1.060531758702533,1.038815350263172,1.045852687705847,1.051239546380595,1.057127187961138,1.063508539089248,1.019095979497993)
}

foreach ($r in $newValues.Keys) {
    $vals = $newValues[$r]
    for ($i = 0; $i -lt $columns.Count; $i++) {
        $ws.Cells.Item($r, $columns[$i]).Value = $vals[$i]
    }
}

Write-Output "Updated vm_pu results ($($newValues.Keys.Count) rows) for the 380 kV case."
